$p = $ppt.ActivePresentation

# The deck's first slide (the "Less developed worlds" Johari-window table, complete
# with its four explanatory annotations) is duplicated, and the duplicate is placed
# immediately after the original as the new slide #2. PowerPoint's default
# Slide.Duplicate() behavior inserts the copy right after the source slide, which
# pushes the former slide #2 ("Less developed worlds" without annotations) and
# slide #3 ("Developed worlds") down to positions #3 and #4 respectively - exactly
# matching the target slide order. No other content changes are required.
$s = $p.Slides.Item(1)
[void]$s.Duplicate()

Write-Output ("Slides.Count=" + $p.Slides.Count)
